$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$textCells = @("D5","D6","D9","D10","D11","D12","D14","D15","D18","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D34","D35","D37","D38","D39","D41","D45","D46","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.800.50'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '2.296.00'
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '299.91'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").Value = '96.76'
$ws.Range("E6").Value = '  -4.69%  '
$ws.Range("E7").Value = '  -1.26%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  -4.00%  '
$ws.Range("D10").Value = '33.33'
$ws.Range("E10").Value = '  -5.18%  '
$ws.Range("D11").Value = '0.0796'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '49.18'
$ws.Range("E12").Value = '  -4.99%  '
$ws.Range("E13").Value = '  +2.25%  '
$ws.Range("D14").Value = '16.74'
$ws.Range("E14").Value = '  +7.20%  '
$ws.Range("D15").Value = '6.76'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").Value = '2.654.41'
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("D17").Value = '2.297.54'
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("D18").Value = '0.805'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = '42.756.70'
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").Value = '0.0₃0899'
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("D21").Value = '11.52'
$ws.Range("E21").Value = '  -2.12%  '
$ws.Range("D22").Value = '6.01'
$ws.Range("E22").Value = '  -1.75%  '
$ws.Range("D23").Value = '67.17'
$ws.Range("E23").Value = '  -1.12%  '
$ws.Range("D24").Value = '235.79'
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("D25").Value = '2.00'
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -3.41%  '
$ws.Range("D28").Value = '24.30'
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").Value = '167.41'
$ws.Range("E29").Value = '  +1.45%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.05'
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '33.73'
$ws.Range("E31").Value = '  -2.46%  '
$ws.Range("D32").Value = '9.08'
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '4.69'
$ws.Range("E34").Value = '  +3.72%  '
$ws.Range("D35").Value = '4.92'
$ws.Range("E35").Value = '  -2.82%  '
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").Value = '16.72'
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").Value = '0.0690'
$ws.Range("E38").Value = '  -2.37%  '
$ws.Range("D39").Value = '2.80'
$ws.Range("E39").Value = '  -3.51%  '
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("D41").Value = '1.74'
$ws.Range("E41").Value = '  -4.63%  '
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("E43").Value = '  -2.79%  '
$ws.Range("D44").Value = '1.985.75'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '0.0279'
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("D46").Value = '9.78'
$ws.Range("E46").Value = '  -1.41%  '
$ws.Range("D47").Value = '17.47'
$ws.Range("E47").Value = '  -5.90%  '
$ws.Range("E48").Value = '  -3.71%  '
$ws.Range("D49").Value = '2.522.74'
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("D50").Value = '52.62'
$ws.Range("E50").Value = '  -6.07%  '
$ws.Range("D51").Value = '4.55'
$ws.Range("E51").Value = '  -7.38%  '
